$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 5 new item rows right after the existing item row (row 7) ---
# This pushes the old "total" row (8) and "footer" row (9) down to rows 13 and 14.
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows.Item(8).Insert()
}

# Clone the formatting/layout of the first item row (row 7) into the five new rows (8-12)
for ($r = 8; $r -le 12; $r++) {
    $ws.Range("A7:Q7").Copy($ws.Range("A" + $r + ":Q" + $r))
}

# Restore the correct row heights for the newly inserted rows
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 25.5

# Helper: assign a value as TEXT, regardless of the cell's underlying number format
function Set-TextValue($range, $value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# --- Row 8: GOLDEN FER 10 MG/ML SYRUP 100 ML ---
$ws.Range("A8").Value = 2
Set-TextValue $ws.Range("C8") "GOLDEN FER 10 MG/ML SYRUP 100 ML"
Set-TextValue $ws.Range("H8") "0:0"
Set-TextValue $ws.Range("L8") "1"
Set-TextValue $ws.Range("N8") "34.00"
Set-TextValue $ws.Range("P8") "34.0000"
Set-TextValue $ws.Range("Q8") "1:0"

# --- Row 9: OMEGA ZAD SYURP 100 ML ---
$ws.Range("A9").Value = 3
Set-TextValue $ws.Range("C9") "OMEGA ZAD SYURP 100 ML"
Set-TextValue $ws.Range("H9") "0:0"
Set-TextValue $ws.Range("L9") "1"
Set-TextValue $ws.Range("N9") "140.00"
Set-TextValue $ws.Range("P9") "140.0000"
Set-TextValue $ws.Range("Q9") "1:0"

# --- Row 10: OXALEPTAL 60MG/ML ORAL SUSP. 100ML ---
$ws.Range("A10").Value = 4
Set-TextValue $ws.Range("C10") "OXALEPTAL 60MG/ML ORAL SUSP. 100ML"
Set-TextValue $ws.Range("H10") "1:0"
Set-TextValue $ws.Range("L10") "1"
Set-TextValue $ws.Range("N10") "89.00"
Set-TextValue $ws.Range("P10") "89.0000"
Set-TextValue $ws.Range("Q10") "1:0"

# --- Row 11: TIRATAM 100MG/ML ORAL SOLUTION 120 ML ---
$ws.Range("A11").Value = 5
Set-TextValue $ws.Range("C11") "TIRATAM 100MG/ML ORAL SOLUTION 120 ML"
Set-TextValue $ws.Range("H11") "1:0"
Set-TextValue $ws.Range("L11") "1"
Set-TextValue $ws.Range("N11") "120.00"
Set-TextValue $ws.Range("P11") "120.0000"
Set-TextValue $ws.Range("Q11") "1:0"

# --- Row 12: VIDROP 2800 I.U./ML ORAL DROPS 15 ML ---
$ws.Range("A12").Value = 6
Set-TextValue $ws.Range("C12") "VIDROP 2800 I.U./ML ORAL DROPS 15 ML"
Set-TextValue $ws.Range("H12") "9:0"
Set-TextValue $ws.Range("L12") "1"
Set-TextValue $ws.Range("N12") "26.00"
Set-TextValue $ws.Range("P12") "26.0000"
Set-TextValue $ws.Range("Q12") "1:0"

# --- Row 13: grand total (shifted down from row 8) ---
$ws.Range("P13").Value = 425.82999999999998

# --- Row 14: footer (shifted down from row 9) - refresh the generated timestamp ---
$ws.Range("A14").Value = "Monday, 29 September, 2025 9:33 AM"
